# quarterly.xlsx update: add 5 earlier quarters (1399/06 .. 1400/06) ahead of
# the existing 5 quarters, shifting the old data from columns E:I to J:N, and
# populate the freshly-inserted E:I columns with the newly reported figures
# ("update database and change read_price algorithm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 new columns right before the first quarter column (E).
#    This shifts the existing quarter columns (E:I) to (J:N) and carries
#    their formatting/styles along, leaving new blank E:I cells that inherit
#    the same per-row styles.
$ws.Range("E1:I1").EntireColumn.Insert()

# 2) New quarter header labels for the freshly inserted columns, on both
#    header rows (row 8 and row 24).
$ws.Range("E8").Value  = "فصل دوم منتهی به 1399/06"
$ws.Range("F8").Value  = "فصل سوم منتهی به 1399/09"
$ws.Range("G8").Value  = "فصل چهارم منتهی به 1399/12"
$ws.Range("H8").Value  = "فصل اول منتهی به 1400/03"
$ws.Range("I8").Value  = "فصل دوم منتهی به 1400/06"

$ws.Range("E24").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("F24").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("G24").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("H24").Value = "فصل اول منتهی به 1400/03"
$ws.Range("I24").Value = "فصل دوم منتهی به 1400/06"

# 3) New figures for the newly inserted quarter columns E:I, row by row.

# هزینه حمل و نقل و انتقال
$ws.Range("E10").Value = 46055
$ws.Range("F10").Value = 124162
$ws.Range("G10").Value = 305040
$ws.Range("H10").Value = 236694
$ws.Range("I10").Value = 178013

# هزینه خدمات پس از فروش
$ws.Range("E11:I11").Value = 0

# حق العمل و کمیسیون فروش
$ws.Range("E12:I12").Value = 0

# هزینه تبلیغات
$ws.Range("E13:I13").Value = 0

# هزینه مواد مصرفی
$ws.Range("E14:I14").Value = 0

# هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Range("E15:I15").Value = 0

# هزینه استهلاک
$ws.Range("E16").Value = 3202
$ws.Range("F16").Value = 6221
$ws.Range("G16").Value = 8776
$ws.Range("H16").Value = 5296
$ws.Range("I16").Value = 6877

# هزینه حقوق و دستمزد
$ws.Range("E17").Value = 31456
$ws.Range("F17").Value = 46173
$ws.Range("G17").Value = 96168
$ws.Range("H17").Value = 56542
$ws.Range("I17").Value = 78485

# هزینه مطالبات مشکوک الوصول
$ws.Range("E18:I18").Value = 0

# سایر هزینه ها
$ws.Range("E19").Value = 322819
$ws.Range("F19").Value = 384141
$ws.Range("G19").Value = 459919
$ws.Range("H19").Value = 449016
$ws.Range("I19").Value = 261277

# جمع
$ws.Range("E20").Value = 403532
$ws.Range("F20").Value = 560697
$ws.Range("G20").Value = 869903
$ws.Range("H20").Value = 747548
$ws.Range("I20").Value = 524652

# تعداد پرسنل غیر تولیدی شرکت
$ws.Range("E26").Value = 402
$ws.Range("F26").Value = 402
$ws.Range("G26").Value = 402
$ws.Range("H26").Value = 402
$ws.Range("I26").Value = 402

# تعداد پرسنل تولیدی شرکت
$ws.Range("E27").Value = 1343
$ws.Range("F27").Value = 1343
$ws.Range("G27").Value = 1358
$ws.Range("H27").Value = 1522
$ws.Range("I27").Value = 1540
